$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps exact text formatting (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "51.532.15"
$ws.Range("E2").Value = "  -0.93%  "

# Row 3
$ws.Range("D3").Value = "2.937.99"
$ws.Range("E3").Value = "  -2.21%  "

# Row 4
$ws.Range("E4").Value = "  -0.21%  "

# Row 5
$ws.Range("D5").Value = "376.48"
$ws.Range("E5").Value = "  +6.03%  "

# Row 6
$ws.Range("D6").Value = "104.34"
$ws.Range("E6").Value = "  -2.24%  "

# Row 7
$ws.Range("D7").Value = "0.542"
$ws.Range("E7").Value = "  -2.48%  "

# Row 8
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.27%  "

# Row 9
$ws.Range("D9").Value = "0.587"
$ws.Range("E9").Value = "  -3.31%  "

# Row 10
$ws.Range("D10").Value = "37.01"
$ws.Range("E10").Value = "  -2.63%  "

# Row 11
$ws.Range("E11").Value = "  -0.79%  "

# Row 12
$ws.Range("D12").Value = "0.0838"
$ws.Range("E12").Value = "  -1.84%  "

# Row 13
$ws.Range("D13").Value = "18.43"
$ws.Range("E13").Value = "  -3.04%  "

# Row 14
$ws.Range("D14").Value = "3.392.44"
$ws.Range("E14").Value = "  -2.48%  "

# Row 15
$ws.Range("D15").Value = "7.40"
$ws.Range("E15").Value = "  -2.39%  "

# Row 16
$ws.Range("D16").Value = "2.930.17"
$ws.Range("E16").Value = "  -1.87%  "

# Row 17
$ws.Range("D17").Value = "0.945"
$ws.Range("E17").Value = "  -5.74%  "

# Row 18
$ws.Range("D18").Value = "51.442.42"
$ws.Range("E18").Value = "  -1.12%  "

# Row 19
$ws.Range("E19").Value = "  +1.20%  "

# Row 20
$ws.Range("D20").Value = "7.34"
$ws.Range("E20").Value = "  -1.44%  "

# Row 21
$ws.Range("D21").Value = "13.03"
$ws.Range("E21").Value = "  -3.36%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0950"
$ws.Range("E22").Value = "  -1.92%  "

# Row 23
$ws.Range("D23").Value = "68.40"
$ws.Range("E23").Value = "  -0.99%  "

# Row 24
$ws.Range("D24").Value = "262.16"
$ws.Range("E24").Value = "  -0.60%  "

# Row 25
$ws.Range("D25").Value = "2.82"
$ws.Range("E25").Value = "  +3.97%  "

# Row 26
$ws.Range("E26").Value = "  -4.84%  "

# Row 27
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "7.17"
$ws.Range("E27").Value = "  +12.99%  "

# Row 28
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.168"
$ws.Range("E28").Value = "  -6.07%  "

# Row 29
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.01%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "25.87"
$ws.Range("E30").Value = "  -3.89%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "7.35"
$ws.Range("E31").Value = "  -0.19%  "

# Row 32
$ws.Range("D32").Value = "0.102"
$ws.Range("E32").Value = "  -6.71%  "

# Row 33
$ws.Range("D33").Value = "9.84"
$ws.Range("E33").Value = "  -3.02%  "

# Row 34
$ws.Range("D34").Value = "51.98"
$ws.Range("E34").Value = "  +1.43%  "

# Row 35
$ws.Range("E35").Value = "  -3.49%  "

# Row 36
$ws.Range("D36").Value = "34.13"
$ws.Range("E36").Value = "  -5.45%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0430"
$ws.Range("E37").Value = "  -0.63%  "

# Row 38
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.29%  "

# Row 39
$ws.Range("D39").Value = "3.03"
$ws.Range("E39").Value = "  -7.65%  "

# Row 40
$ws.Range("D40").Value = "17.03"
$ws.Range("E40").Value = "  -2.39%  "

# Row 41
$ws.Range("D41").Value = "2.61"
$ws.Range("E41").Value = "  -8.58%  "

# Row 42
$ws.Range("D42").Value = "1.83"
$ws.Range("E42").Value = "  -5.57%  "

# Row 43
$ws.Range("D43").Value = "0.115"
$ws.Range("E43").Value = "  -1.65%  "

# Row 44
$ws.Range("D44").Value = "124.46"
$ws.Range("E44").Value = "  +0.54%  "

# Row 45
$ws.Range("D45").Value = "21.95"
$ws.Range("E45").Value = "  -5.52%  "

# Row 46
$ws.Range("D46").Value = "2.06"
$ws.Range("E46").Value = "  -4.96%  "

# Row 47
$ws.Range("D47").Value = "0.272"
$ws.Range("E47").Value = "  +12.03%  "

# Row 48
$ws.Range("D48").Value = "2.024.63"
$ws.Range("E48").Value = "  -4.73%  "

# Row 49
$ws.Range("E49").Value = "  -0.83%  "

# Row 50
$ws.Range("D50").Value = "3.19"
$ws.Range("E50").Value = "  -3.62%  "

# Row 51
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "0.0324"
$ws.Range("E51").Value = "  -2.36%  "
